$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1080.8413
$ws.Range("I15").Value = 1080.8413
$ws.Range("K15").Value = 3242.5239
$ws.Range("M15").Value = -3073.5239
$ws.Range("H40").Value = 83335960
$ws.Range("I40").Value = 2633.3333
$ws.Range("K40").Value = 2633.3333
$ws.Range("M40").Value = -2458.3333
$ws.Range("H69").Value = 14285.214
$ws.Range("J69").Value = 14285.214
$ws.Range("L69").Value = 42855.642
$ws.Range("N69").Value = -44603.642
$ws.Range("H72").Value = 14285.214
$ws.Range("J72").Value = 14285.214
$ws.Range("L72").Value = 128566.926
$ws.Range("N72").Value = -137302.926
$ws.Range("H88").Value = 1535.2593
$ws.Range("I88").Value = 1278.4166
$ws.Range("J88").Value = 1740.7333
$ws.Range("K88").Value = 1278.4166
$ws.Range("L88").Value = 1740.7333
$ws.Range("M88").Value = -872.4166
$ws.Range("N88").Value = -2552.7333
$ws.Range("H91").Value = 1535.2593
$ws.Range("I91").Value = 1278.4166
$ws.Range("J91").Value = 1740.7333
$ws.Range("K91").Value = 1278.4166
$ws.Range("L91").Value = 1740.7333
$ws.Range("M91").Value = 125.5834
$ws.Range("N91").Value = -4548.7333
$ws.Range("H98").Value = 2156.7778
$ws.Range("I98").Value = 1707.75
$ws.Range("K98").Value = 1707.75
$ws.Range("M98").Value = -209.75
$ws.Range("H111").Value = 11099.5
$ws.Range("J111").Value = 11099.5
$ws.Range("L111").Value = 33298.5
$ws.Range("N111").Value = -39432.5
$ws.Range("H122").Value = 2156.7778
$ws.Range("I122").Value = 1707.75
$ws.Range("K122").Value = 5123.25
$ws.Range("M122").Value = -2673.25
$ws.Range("H132").Value = 2028.4048
$ws.Range("I132").Value = 2030.1
$ws.Range("K132").Value = 6090.299999999999
$ws.Range("M132").Value = -3560.299999999999
$ws.Range("H137").Value = 2249.0908
$ws.Range("I137").Value = 2373.7
$ws.Range("J137").Value = 1003
$ws.Range("K137").Value = 7121.099999999999
$ws.Range("L137").Value = 3009
$ws.Range("M137").Value = -4571.099999999999
$ws.Range("N137").Value = -8109
$ws.Range("H138").Value = 5585.077
$ws.Range("I138").Value = 4567.2856
$ws.Range("K138").Value = 13701.8568
$ws.Range("M138").Value = -8561.856800000001
$ws.Range("H141").Value = 3421.5715
$ws.Range("I141").Value = 3592.077
$ws.Range("K141").Value = 10776.231
$ws.Range("M141").Value = -5596.231

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8735.344999999999
$ws.Range("I32").Value = 8341.161
$ws.Range("K32").Value = 8341.161
$ws.Range("M32").Value = -8054.161
$ws.Range("H45").Value = 2881
$ws.Range("I45").Value = 1171.3334
$ws.Range("K45").Value = 1171.3334
$ws.Range("M45").Value = -794.3334
$ws.Range("H61").Value = 4085518.8
$ws.Range("I61").Value = 4881620
$ws.Range("K61").Value = 4881620
$ws.Range("M61").Value = -4881408
$ws.Range("H110").Value = 5983.409
$ws.Range("I110").Value = 6043.6875
$ws.Range("K110").Value = 6043.6875
$ws.Range("M110").Value = -3998.6875
$ws.Range("H132").Value = 2809.611
$ws.Range("I132").Value = 2741.7188
$ws.Range("K132").Value = 8225.1564
$ws.Range("M132").Value = -5695.1564
$ws.Range("H136").Value = 4085518.8
$ws.Range("I136").Value = 4881620
$ws.Range("K136").Value = 14644860
$ws.Range("M136").Value = -14642310
$ws.Range("H141").Value = 115999
$ws.Range("I141").Value = 112000
$ws.Range("K141").Value = 112000
$ws.Range("M141").Value = -106820

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H59").Value = 90220
$ws.Range("I59").Value = 55440.5
$ws.Range("J59").Value = 124999.5
$ws.Range("K59").Value = 55440.5
$ws.Range("L59").Value = 124999.5
$ws.Range("M59").Value = -54295.5
$ws.Range("N59").Value = -127289.5
$ws.Range("H68").Value = 75000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 75000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H99").Value = 14808.714
$ws.Range("I99").Value = 7403.5293
$ws.Range("K99").Value = 7403.5293
$ws.Range("M99").Value = -5905.5293
$ws.Range("H107").Value = 2168.3333
$ws.Range("I107").Value = 756.6667
$ws.Range("J107").Value = 3580
$ws.Range("K107").Value = 756.6667
$ws.Range("L107").Value = 3580
$ws.Range("M107").Value = 1163.3333
$ws.Range("N107").Value = -7420
$ws.Range("H126").Value = 14808.714
$ws.Range("I126").Value = 7403.5293
$ws.Range("K126").Value = 22210.5879
$ws.Range("M126").Value = -19740.5879
$ws.Range("H132").Value = 1662.8572
$ws.Range("I132").Value = 1483.3182
$ws.Range("K132").Value = 4449.9546
$ws.Range("M132").Value = -1919.9546

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 22222450
$ws.Range("I23").Value = 190.4
$ws.Range("J23").Value = 33333582
$ws.Range("K23").Value = 571.2
$ws.Range("L23").Value = 100000746
$ws.Range("M23").Value = -336.2
$ws.Range("N23").Value = -100001216
$ws.Range("H113").Value = 1798.9412
$ws.Range("I113").Value = 2305.75
$ws.Range("J113").Value = 1643
$ws.Range("K113").Value = 6917.25
$ws.Range("L113").Value = 4929
$ws.Range("M113").Value = -4747.25
$ws.Range("N113").Value = -9269
$ws.Range("H137").Value = 7510.756
$ws.Range("I137").Value = 4617.5
$ws.Range("J137").Value = 9362.440000000001
$ws.Range("K137").Value = 13852.5
$ws.Range("L137").Value = 28087.32
$ws.Range("M137").Value = -8752.5
$ws.Range("N137").Value = -38287.32

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1519.0834
$ws.Range("I97").Value = 1454.7778
$ws.Range("K97").Value = 1454.7778
$ws.Range("M97").Value = -958.7778000000001
$ws.Range("H132").Value = 2892.0938
$ws.Range("I132").Value = 3166.2
$ws.Range("J132").Value = 1913.1428
$ws.Range("K132").Value = 9498.599999999999
$ws.Range("L132").Value = 5739.428400000001
$ws.Range("M132").Value = -6968.599999999999
$ws.Range("N132").Value = -10799.4284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 18858570
$ws.Range("I22").Value = 26400318
$ws.Range("J22").Value = 4200
$ws.Range("K22").Value = 26400318
$ws.Range("L22").Value = 4200
$ws.Range("M22").Value = -26400023
$ws.Range("N22").Value = -4790
$ws.Range("H27").Value = 18858570
$ws.Range("I27").Value = 26400318
$ws.Range("J27").Value = 4200
$ws.Range("K27").Value = 26400318
$ws.Range("L27").Value = 4200
$ws.Range("M27").Value = -26400211
$ws.Range("N27").Value = -4414
$ws.Range("H46").Value = 1466.25
$ws.Range("I46").Value = 1297.6
$ws.Range("J46").Value = 1747.3334
$ws.Range("K46").Value = 1297.6
$ws.Range("L46").Value = 1747.3334
$ws.Range("M46").Value = -1109.6
$ws.Range("N46").Value = -2123.3334
$ws.Range("H133").Value = 124000
$ws.Range("J133").Value = 124000
$ws.Range("L133").Value = 124000
$ws.Range("N133").Value = -129060

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 99996.5
$ws.Range("J16").Value = 99996.5
$ws.Range("L16").Value = 99996.5
$ws.Range("N16").Value = -100580.5
$ws.Range("H63").Value = 39250
$ws.Range("J63").Value = 39250
$ws.Range("L63").Value = 39250
$ws.Range("N63").Value = -40498
$ws.Range("H66").Value = 39250
$ws.Range("J66").Value = 39250
$ws.Range("L66").Value = 117750
$ws.Range("N66").Value = -123990
$ws.Range("H81").Value = 3567.8
$ws.Range("I81").Value = 2914.2222
$ws.Range("J81").Value = 9450
$ws.Range("K81").Value = 5828.4444
$ws.Range("L81").Value = 18900
$ws.Range("M81").Value = -4767.4444
$ws.Range("N81").Value = -21022
$ws.Range("H84").Value = 3567.8
$ws.Range("I84").Value = 2914.2222
$ws.Range("J84").Value = 9450
$ws.Range("K84").Value = 29142.222
$ws.Range("L84").Value = 94500
$ws.Range("M84").Value = -23838.222
$ws.Range("N84").Value = -105108
$ws.Range("H132").Value = 3273.92
$ws.Range("I132").Value = 2040.7142
$ws.Range("J132").Value = 9748.25
$ws.Range("K132").Value = 6122.142599999999
$ws.Range("L132").Value = 29244.75
$ws.Range("M132").Value = -3592.142599999999
$ws.Range("N132").Value = -34304.75
